$d = $word.ActiveDocument

# The paragraph that currently reads:
#   "O Hipermercado Edmélio precisa registrar suas vendas! A super inauguração
#    se aproxima ... para registrar os produt[_GoBack]os o mais rápido possível."
# gets split into two paragraphs: an empty leading paragraph that keeps the
# original paragraph identity/rsids and now only holds the (moved) _GoBack
# bookmark, followed by a brand new paragraph carrying all of the original
# text (with the run that used to be split by the bookmark now merged back
# into a single run).

$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "O Hiper*possível.*") {
        $target = $p
        break
    }
}

$full = $d.Range($target.Range.Start, $target.Range.End - 1)

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
'<w:body>' +
'<w:p w:rsidR="00DA66C1" w:rsidRDefault="006836AA" w:rsidP="002626D7"><w:pPr><w:jc w:val="both"/></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>' +
'<w:p><w:pPr><w:jc w:val="both"/></w:pPr>' +
'<w:r><w:t>O Hiper</w:t></w:r>' +
'<w:r><w:t>m</w:t></w:r>' +
'<w:r><w:t xml:space="preserve">ercado </w:t></w:r>' +
'<w:proofErr w:type="spellStart"/>' +
'<w:r><w:t>Edm&#233;lio</w:t></w:r>' +
'<w:proofErr w:type="spellEnd"/>' +
'<w:r><w:t xml:space="preserve"> precisa registrar suas vendas! A </w:t></w:r>' +
'<w:proofErr w:type="spellStart"/>' +
'<w:r><w:t>super</w:t></w:r>' +
'<w:proofErr w:type="spellEnd"/>' +
'<w:r><w:t xml:space="preserve"> inaugura&#231;&#227;o se aproxima e s&#227;o esperados muitos clientes. Estes ser&#227;o atendidos em 5 caixas que contar&#227;o com um operador cada. Os operadores precisam de aux&#237;lio para registrar os produtos o mais r&#225;pido poss&#237;vel.</w:t></w:r>' +
'</w:p>' +
'</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$full.InsertXML($xml)
